$d = $word.ActiveDocument
$r = $d.Content
$found = $r.Find.Execute("submit your scores", $true, $false, $false, $false, $false, $true, 1, $false, "submit your results", 2)
Write-Host "Found: $found"
